# Refactor eavl_fpr_latency: complete function extraction from function
#
# Both worksheets ("fpr_latency_tradeoff" and "fpr_sdr_tradeoff") had a
# single data row (row 2, columns A:F) holding per-class metrics plus the
# target_fpr threshold in column F. The metrics now live on row 3 shifted
# one column to the right (column A now holds what used to be in F, i.e.
# the previous target_fpr value), while row 2 keeps only the new
# target_fpr value (0.0005) in column F.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: fpr_latency_tradeoff -----------------------------------
$ws1 = $wb.Worksheets.Item("fpr_latency_tradeoff")

$s1_oldA2 = $ws1.Range("A2").Value2
$s1_oldB2 = $ws1.Range("B2").Value2
$s1_oldC2 = $ws1.Range("C2").Value2
$s1_oldD2 = $ws1.Range("D2").Value2
$s1_oldE2 = $ws1.Range("E2").Value2
$s1_oldF2 = $ws1.Range("F2").Value2
$s1_oldF3 = $ws1.Range("F3").Value2

$ws1.Range("A2").Value = ""
$ws1.Range("B2").Value = ""
$ws1.Range("C2").Value = ""
$ws1.Range("D2").Value = ""
$ws1.Range("E2").Value = ""
$ws1.Range("F2").Value = $s1_oldF3

$ws1.Range("A3").Value = $s1_oldF2
$ws1.Range("B3").Value = $s1_oldA2
$ws1.Range("C3").Value = $s1_oldB2
$ws1.Range("D3").Value = $s1_oldC2
$ws1.Range("E3").Value = $s1_oldD2
$ws1.Range("F3").Value = $s1_oldE2

# --- Sheet 2: fpr_sdr_tradeoff ----------------------------------------
$ws2 = $wb.Worksheets.Item("fpr_sdr_tradeoff")

$s2_oldA2 = $ws2.Range("A2").Value2
$s2_oldB2 = $ws2.Range("B2").Value2
$s2_oldC2 = $ws2.Range("C2").Value2
$s2_oldD2 = $ws2.Range("D2").Value2
$s2_oldE2 = $ws2.Range("E2").Value2
$s2_oldF3 = $ws2.Range("F3").Value2

$ws2.Range("A2").Value = ""
$ws2.Range("B2").Value = ""
$ws2.Range("C2").Value = ""
$ws2.Range("D2").Value = ""
$ws2.Range("E2").Value = ""
$ws2.Range("F2").Value = $s2_oldF3

$ws2.Range("A3").Value = $s2_oldA2
$ws2.Range("B3").Value = $s2_oldB2
$ws2.Range("C3").Value = $s2_oldC2
$ws2.Range("D3").Value = $s2_oldD2
$ws2.Range("E3").Value = $s2_oldE2
$ws2.Range("F3").Value = $s1_oldE2
